$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'258.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.33%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.79%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.800"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-10.45%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05970"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.53%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.689"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.33%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8755"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.44%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9518"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.51%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1415"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.58%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.03612"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'5.27%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.08%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'-1.23%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09237"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.001547"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.05%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006068"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.35%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005979"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.27%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.485"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.33%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.229"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.11%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.219"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.52%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3134"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.02%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1290"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.06%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.528"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.05%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04222"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.55%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").Value = "'0.001222"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.42%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004513"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-11.97%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-0.04%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'-22.97%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03845"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.08%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006018"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.63%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1102"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.17%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002199"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-7.99%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01080"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'9.15%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005492"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.86%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'8.97%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002125"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-3.75%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("E50").Style = "Normal"
